$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying source data for this weekly Fruta/Hortaliza sheet was
# re-synced, which re-orders the date-keyed rows (2-19). Row 12 is unaffected.
# Apply the new values cell by cell, column by column, to match the refreshed data.

# Row 2 (was row 6 data in the prior export)
$ws.Cells.Item(2, "D").Value = 44628
$ws.Cells.Item(2, "M").Value = 40
$ws.Cells.Item(2, "N").Value = 6000
$ws.Cells.Item(2, "O").Value = 6000
$ws.Cells.Item(2, "P").Value = 6000
$ws.Cells.Item(2, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(2, "S").Value = 3000

# Row 3 (was row 17 data in the prior export)
$ws.Cells.Item(3, "D").Value = 44627
$ws.Cells.Item(3, "M").Value = 45
$ws.Cells.Item(3, "N").Value = 6000
$ws.Cells.Item(3, "O").Value = 6000
$ws.Cells.Item(3, "P").Value = 6000
$ws.Cells.Item(3, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(3, "S").Value = 3000

# Row 4 (was row 7 data in the prior export)
$ws.Cells.Item(4, "D").Value = 44587
$ws.Cells.Item(4, "M").Value = 165
$ws.Cells.Item(4, "N").Value = 6500
$ws.Cells.Item(4, "O").Value = 7000
$ws.Cells.Item(4, "P").Value = 6742
$ws.Cells.Item(4, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(4, "S").Value = 3371

# Row 5 (was row 19 data in the prior export)
$ws.Cells.Item(5, "D").Value = 44214
$ws.Cells.Item(5, "M").Value = 48
$ws.Cells.Item(5, "N").Value = 6000
$ws.Cells.Item(5, "O").Value = 6000
$ws.Cells.Item(5, "P").Value = 6000
$ws.Cells.Item(5, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(5, "S").Value = 3000

# Row 6 (was row 2 data in the prior export)
$ws.Cells.Item(6, "D").Value = 44211
$ws.Cells.Item(6, "M").Value = 45
$ws.Cells.Item(6, "N").Value = 6000
$ws.Cells.Item(6, "O").Value = 6000
$ws.Cells.Item(6, "P").Value = 6000
$ws.Cells.Item(6, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(6, "S").Value = 3000

# Row 7 (was row 14 data in the prior export)
$ws.Cells.Item(7, "D").Value = 44582
$ws.Cells.Item(7, "M").Value = 150
$ws.Cells.Item(7, "N").Value = 6000
$ws.Cells.Item(7, "O").Value = 6500
$ws.Cells.Item(7, "P").Value = 6233
$ws.Cells.Item(7, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(7, "S").Value = 3116

# Row 8 (was row 11 data in the prior export)
$ws.Cells.Item(8, "D").Value = 44209
$ws.Cells.Item(8, "M").Value = 58
$ws.Cells.Item(8, "N").Value = 6000
$ws.Cells.Item(8, "O").Value = 6000
$ws.Cells.Item(8, "P").Value = 6000
$ws.Cells.Item(8, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(8, "S").Value = 3000

# Row 9 (was row 18 data in the prior export)
$ws.Cells.Item(9, "D").Value = 44592
$ws.Cells.Item(9, "M").Value = 30
$ws.Cells.Item(9, "N").Value = 8000
$ws.Cells.Item(9, "O").Value = 8000
$ws.Cells.Item(9, "P").Value = 8000
$ws.Cells.Item(9, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(9, "S").Value = 4000

# Row 10 (was row 9 data in the prior export)
$ws.Cells.Item(10, "D").Value = 44960
$ws.Cells.Item(10, "M").Value = 40
$ws.Cells.Item(10, "N").Value = 7000
$ws.Cells.Item(10, "O").Value = 7000
$ws.Cells.Item(10, "P").Value = 7000
$ws.Cells.Item(10, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(10, "S").Value = 3500

# Row 11 (was row 13 data in the prior export)
$ws.Cells.Item(11, "D").Value = 44606
$ws.Cells.Item(11, "M").Value = 45
$ws.Cells.Item(11, "N").Value = 7000
$ws.Cells.Item(11, "O").Value = 7000
$ws.Cells.Item(11, "P").Value = 7000
$ws.Cells.Item(11, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(11, "S").Value = 3500

# Row 13 (was row 16 data in the prior export)
$ws.Cells.Item(13, "D").Value = 44614
$ws.Cells.Item(13, "M").Value = 45
$ws.Cells.Item(13, "N").Value = 6000
$ws.Cells.Item(13, "O").Value = 6000
$ws.Cells.Item(13, "P").Value = 6000
$ws.Cells.Item(13, "R").Value = 'Provincia de Linares'
$ws.Cells.Item(13, "S").Value = 3000

# Row 14 (was row 8 data in the prior export)
$ws.Cells.Item(14, "D").Value = 44959
$ws.Cells.Item(14, "M").Value = 40
$ws.Cells.Item(14, "N").Value = 7000
$ws.Cells.Item(14, "O").Value = 7000
$ws.Cells.Item(14, "P").Value = 7000
$ws.Cells.Item(14, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(14, "S").Value = 3500

# Row 15 (was row 10 data in the prior export)
$ws.Cells.Item(15, "D").Value = 45001
$ws.Cells.Item(15, "M").Value = 66
$ws.Cells.Item(15, "N").Value = 7500
$ws.Cells.Item(15, "O").Value = 8000
$ws.Cells.Item(15, "P").Value = 7773
$ws.Cells.Item(15, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(15, "S").Value = 3886

# Row 16 (was row 3 data in the prior export)
$ws.Cells.Item(16, "D").Value = 44974
$ws.Cells.Item(16, "M").Value = 130
$ws.Cells.Item(16, "N").Value = 7000
$ws.Cells.Item(16, "O").Value = 7500
$ws.Cells.Item(16, "P").Value = 7269
$ws.Cells.Item(16, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(16, "S").Value = 3634

# Row 17 (was row 4 data in the prior export)
$ws.Cells.Item(17, "D").Value = 44585
$ws.Cells.Item(17, "M").Value = 160
$ws.Cells.Item(17, "N").Value = 6500
$ws.Cells.Item(17, "O").Value = 7000
$ws.Cells.Item(17, "P").Value = 6750
$ws.Cells.Item(17, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(17, "S").Value = 3375

# Row 18 (was row 5 data in the prior export)
$ws.Cells.Item(18, "D").Value = 44589
$ws.Cells.Item(18, "M").Value = 60
$ws.Cells.Item(18, "N").Value = 6000
$ws.Cells.Item(18, "O").Value = 6000
$ws.Cells.Item(18, "P").Value = 6000
$ws.Cells.Item(18, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(18, "S").Value = 3000

# Row 19 (was row 15 data in the prior export)
$ws.Cells.Item(19, "D").Value = 44588
$ws.Cells.Item(19, "M").Value = 160
$ws.Cells.Item(19, "N").Value = 6500
$ws.Cells.Item(19, "O").Value = 7000
$ws.Cells.Item(19, "P").Value = 6750
$ws.Cells.Item(19, "R").Value = 'Provincia de Curicó'
$ws.Cells.Item(19, "S").Value = 3375
